$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = 1.5
$ws.Range("M3").Value = 2.63
$ws.Range("T3").Value = 6.5
$ws.Range("U3").Value = 11
$ws.Range("AB5").Value = 17
$ws.Range("AE5").Value = 15
$ws.Range("AI5").Value = 41
$ws.Range("G5").Value = 1.62
$ws.Range("H5").Value = 4.2
$ws.Range("U5").Value = 8
$ws.Range("W5").Value = 12
$ws.Range("AA6").Value = 7.9
$ws.Range("AB6").Value = 12.5
$ws.Range("AE6").Value = 19
$ws.Range("AG6").Value = 16
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 4.3
$ws.Range("T6").Value = 8
$ws.Range("V6").Value = 7
$ws.Range("W6").Value = 8.75
$ws.Range("Y6").Value = 15.5
$ws.Range("Z6").Value = 16.5
$ws.Range("AB7").Value = 15.5
$ws.Range("AC7").Value = 60
$ws.Range("AD7").Value = 350
$ws.Range("AE7").Value = 23
$ws.Range("AF7").Value = 65
$ws.Range("AG7").Value = 23
$ws.Range("AJ7").Value = 60
$ws.Range("H7").Value = 4.45
$ws.Range("I7").Value = 9.5
$ws.Range("M7").Value = 4.5
$ws.Range("N7").Value = 1.55
$ws.Range("O7").Value = 2.15
$ws.Range("T7").Value = 6.3
$ws.Range("U7").Value = 5.5
$ws.Range("V7").Value = 7.1
$ws.Range("W7").Value = 6.9
$ws.Range("AH9").Value = 29
$ws.Range("AI9").Value = 21
$ws.Range("AJ9").Value = 26
$ws.Range("G9").Value = 2.15
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 2.8
$ws.Range("N9").Value = 1.7
$ws.Range("O9").Value = 2.1
$ws.Range("R9").Value = 1.62
$ws.Range("S9").Value = 2.2
$ws.Range("T9").Value = 9.5
$ws.Range("U9").Value = 12
$ws.Range("W9").Value = 21
$ws.Range("AI10").Value = 29
$ws.Range("G10").Value = 1.7
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 3.8
$ws.Range("N10").Value = 1.48
$ws.Range("O10").Value = 2.6
$ws.Range("X10").Value = 12
$ws.Range("AE11").Value = 15
$ws.Range("AF11").Value = 23
$ws.Range("AG11").Value = 15
$ws.Range("AI11").Value = 34
$ws.Range("I11").Value = 4
$ws.Range("L11").Value = 1.2
$ws.Range("M11").Value = 4.33
$ws.Range("N11").Value = 1.7
$ws.Range("O11").Value = 2.1
$ws.Range("G12").Value = 2.1
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 3.2
$ws.Range("L12").Value = 1.25
$ws.Range("M12").Value = 3.75
$ws.Range("N12").Value = 1.88
$ws.Range("O12").Value = 1.98
$ws.Range("P12").Value = 1.36
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = 1.67
$ws.Range("S12").Value = 2.1
$ws.Range("W12").Value = 19
$ws.Range("AA13").Value = 5.9
$ws.Range("AB13").Value = 15.5
$ws.Range("AC13").Value = 90
$ws.Range("AD13").Value = 800
$ws.Range("AF13").Value = 16.5
$ws.Range("AJ13").Value = 45
$ws.Range("H13").Value = 2.95
$ws.Range("I13").Value = 3.25
$ws.Range("J13").Value = 1.09
$ws.Range("K13").Value = 6.1
$ws.Range("L13").Value = 1.42
$ws.Range("M13").Value = 2.67
$ws.Range("N13").Value = 2.22
$ws.Range("O13").Value = 1.6
$ws.Range("P13").Value = 1.53
$ws.Range("Q13").Value = 2.35
$ws.Range("R13").Value = 1.91
$ws.Range("S13").Value = 1.8
$ws.Range("T13").Value = 6.4
$ws.Range("U13").Value = 9.75
$ws.Range("V13").Value = 9.25
$ws.Range("X13").Value = 21
$ws.Range("Y13").Value = 35
$ws.Range("Z13").Value = 6.1
$ws.Range("AA15").Value = 5.5
$ws.Range("AB15").Value = 10.75
$ws.Range("AC15").Value = 40
$ws.Range("AD15").Value = 250
$ws.Range("AE15").Value = 10.5
$ws.Range("AF15").Value = 22
$ws.Range("AG15").Value = 11.5
$ws.Range("AH15").Value = 60
$ws.Range("G15").Value = 1.78
$ws.Range("H15").Value = 3.15
$ws.Range("I15").Value = 4.4
$ws.Range("M15").Value = 3.45
$ws.Range("N15").Value = 1.8
$ws.Range("O15").Value = 1.8
$ws.Range("R15").Value = 1.75
$ws.Range("T15").Value = 6.1
$ws.Range("U15").Value = 7.5
$ws.Range("V15").Value = 6.7
$ws.Range("W15").Value = 12.5
$ws.Range("X15").Value = 11.25
$ws.Range("Y15").Value = 18.5
$ws.Range("Z15").Value = 9.25
$ws.Range("AA16").Value = 17
$ws.Range("AB16").Value = 35
$ws.Range("AC16").Value = 150
$ws.Range("AE16").Value = 50
$ws.Range("AG16").Value = 55
$ws.Range("AI16").Value = 250
$ws.Range("AJ16").Value = 150
$ws.Range("G16").Value = 1.08
$ws.Range("H16").Value = 7.8
$ws.Range("I16").Value = 17
$ws.Range("O16").Value = 3.5
$ws.Range("R16").Value = 2.27
$ws.Range("S16").Value = 1.56
$ws.Range("T16").Value = 9.75
$ws.Range("U16").Value = 6.1
$ws.Range("V16").Value = 11
$ws.Range("W16").Value = 5.5
$ws.Range("X16").Value = 9.75
$ws.Range("Y16").Value = 32
$ws.Range("AI17").Value = 26
$ws.Range("AJ17").Value = 34
$ws.Range("G17").Value = 2.3
$ws.Range("I17").Value = 3.2
$ws.Range("K17").Value = 9
$ws.Range("N17").Value = 2.1
$ws.Range("O17").Value = 1.7
$ws.Range("T17").Value = 7.5
$ws.Range("U17").Value = 11
$ws.Range("AA18").Value = 9.5
$ws.Range("AD18").Value = 201
$ws.Range("AF18").Value = 41
$ws.Range("G18").Value = 1.3
$ws.Range("H18").Value = 4.75
$ws.Range("R18").Value = 1.8
$ws.Range("S18").Value = 1.91
$ws.Range("W18").Value = 9
$ws.Range("AC19").Value = 51
$ws.Range("AD19").Value = 201
$ws.Range("AE19").Value = 13
$ws.Range("AJ19").Value = 41
$ws.Range("G19").Value = 1.67
$ws.Range("I19").Value = 4.33
$ws.Range("J19").Value = 1.05
$ws.Range("K19").Value = 8
$ws.Range("L19").Value = 1.25
$ws.Range("M19").Value = 3.6
$ws.Range("N19").Value = 1.8
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = 1.36
$ws.Range("Q19").Value = 2.9
$ws.Range("R19").Value = 1.72
$ws.Range("S19").Value = 2
$ws.Range("T19").Value = 8
$ws.Range("U19").Value = 8.5
$ws.Range("W19").Value = 13
$ws.Range("Z19").Value = 12
$ws.Range("J20").Value = 1.06
$ws.Range("K20").Value = 10
$ws.Range("N20").Value = 2.03
$ws.Range("O20").Value = 1.83
$ws.Range("AA22").Value = 8
$ws.Range("AB22").Value = 20
$ws.Range("AC22").Value = 110
$ws.Range("AD22").Value = 900
$ws.Range("AF22").Value = 45
$ws.Range("AG22").Value = 22
$ws.Range("AH22").Value = 175
$ws.Range("AI22").Value = 90
$ws.Range("AJ22").Value = 80
$ws.Range("G22").Value = 1.44
$ws.Range("H22").Value = 4.1
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 1.06
$ws.Range("K22").Value = 7.5
$ws.Range("L22").Value = 1.28
$ws.Range("M22").Value = 3.35
$ws.Range("N22").Value = 1.82
$ws.Range("O22").Value = 1.88
$ws.Range("P22").Value = 1.4
$ws.Range("Q22").Value = 2.75
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 1.72
$ws.Range("T22").Value = 6.3
$ws.Range("U22").Value = 6.5
$ws.Range("X22").Value = 12
$ws.Range("Y22").Value = 29
$ws.Range("Z22").Value = 7.5
